$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: doctor wants to enter disease details of the patient -> add respose/destination
$ws.Range("E7").Value = "patient profile"
$ws.Range("F7").Value = "doctor"

# Row 8: nurse wants to check patients profiles -> add respose/destination
$ws.Range("E8").Value = "patient profile"
$ws.Range("F8").Value = "nurse"

# Row 10: nurse wants to enter disease details of the patient -> add respose/destination
$ws.Range("E10").Value = "patient profile"
$ws.Range("F10").Value = "nurse"

# Row 12: lab staff wants to add lab reports -> add respose
$ws.Range("E12").Value = "report"

# Row 14: admission officer wants to create a new record for a new visit -> add respose (destination already set)
$ws.Range("E14").Value = "pateint record"
